# Add "Dada" data: six more rows to the Acc sheet, and a brand new
# "Sheet1" worksheet (after "New Customer") holding a second copy of the
# New-123..New-128 customer-address rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Acc sheet: append rows 11-16 (Row_10 .. Row_15)
# ---------------------------------------------------------------------
$wsAcc = $wb.Worksheets.Item("Acc")
$wsAcc.Select()

$wsAcc.Range("A11").Value = "Row_10"
$wsAcc.Range("B11").Value = "Current"

$wsAcc.Range("A12").Value = "Row_11"
$wsAcc.Range("B12").Value = "Current"
$wsAcc.Range("C12").Value = "Banny"

$wsAcc.Range("A13").Value = "Row_12"
$wsAcc.Range("B13").Value = "Current"
$wsAcc.Range("C13").Value = "Rajesh"

$wsAcc.Range("A14").Value = "Row_13"
$wsAcc.Range("B14").Value = "Saving"
$wsAcc.Range("C14").Value = "Charith"

$wsAcc.Range("A15").Value = "Row_14"
$wsAcc.Range("B15").Value = "Saving"

$wsAcc.Range("A16").Value = "Row_15"
$wsAcc.Range("B16").Value = "Saving"

# Mirrors the saved file's selection on the Acc tab after the edit.
$wsAcc.Range("C12:C14").Select()

# ---------------------------------------------------------------------
# 2) New Customer sheet: selection moves to A2:B6
# ---------------------------------------------------------------------
$wsNC = $wb.Worksheets.Item("New Customer")
$wsNC.Select()
$wsNC.Range("A2:B6").Select()

# ---------------------------------------------------------------------
# 3) Brand new "Sheet1" worksheet, placed after "New Customer"
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNew = $wb.Worksheets.Add($null, $lastSheet)

$wsNew.Range("A1").Value = "New123"
$wsNew.Range("B1").Value = "Chennai1"

$wsNew.Range("A2").Value = "New124"
$wsNew.Range("B2").Value = "Hyderabad1"

$wsNew.Range("A3").Value = "New125"
$wsNew.Range("B3").Value = "Bangalore1"

$wsNew.Range("A4").Value = "New127"
$wsNew.Range("B4").Value = "Piler1"

$wsNew.Range("A5").Value = "New128"
$wsNew.Range("B5").Value = "Tirupati1"

$wsNew.Range("A1:B5").Select()
